$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of CRM accuracy data
$ws.Range("A67").Value = 20220315
$ws.Range("B67").Value = 2225.7539999999999
$ws.Range("C67").Value = 2224.4699999999998
$ws.Range("D67").Formula = "=100*(B67-C67)/C67"
$ws.Range("E67").Value = 180
$ws.Range("F67").Value = "CRM OPENED 20220303"

$ws.Range("A68").Value = 20200317
$ws.Range("B68").Value = 2223.0121947257499
$ws.Range("C68").Value = 2224.4699999999998
$ws.Range("D68").Formula = "=100*(B68-C68)/C68"
$ws.Range("E68").Value = 180
$ws.Range("F68").Value = "CRM OPENED 20220303"

# Update the view so the new rows are visible, matching author's scroll position
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("E64").Select()
